$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 736, shifting existing rows 736:835 down to 737:836
# (dimension grows from A1:R835 to A1:R836)
$ws.Rows.Item(736).Insert()

# Populate the newly inserted row 736 with the new weekly data record
$ws.Range("A736").Value = 11
$ws.Range("B736").Value = 'Vega Monumental Concepción'
$ws.Range("C736").Value = 'Bíobío'
$ws.Range("D736").Value = 45142
$ws.Range("E736").Value = 8
$ws.Range("F736").Value = 100112020
$ws.Range("G736").Value = 'Tomate'
$ws.Range("H736").Value = 'Larga vida'
$ws.Range("I736").Value = 'Primera'
$ws.Range("J736").Value = 450
$ws.Range("K736").Value = 21000
$ws.Range("L736").Value = 22000
$ws.Range("M736").Value = 21444
$ws.Range("N736").Value = '$/bandeja 18 kilos'
$ws.Range("O736").Value = 'Región de Arica y Parinacota'
$ws.Range("P736").Value = 1191
$ws.Range("Q736").Value = 18
$ws.Range("R736").Value = 'Hortaliza'
